# Generate Report for Handoff
# The handoff transform failed for 821d87f1-8817-4972-a132-099f001c662d.md;
# a new markdown file (7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.md) is generated
# instead, and the stale handoff-in-progress data (xlf hand-off files /
# timestamps / "Include" reason) is replaced with "Ignored" / default
# datetime values since the handoff never produced a target file.

$wb = $excel.ActiveWorkbook

$newFileName = "7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.md"
$cfgAddr = "https://github.com/OpenLocalizationTest/oltest/blob/dbda4a81e1d098bcdc75ea17448fd19d95477593/.localization-config"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

$wsOverview.Hyperlinks.Delete()
$mdAddrOverview = "https://github.com/OpenLocalizationTest/oltest/blob/dbda4a81e1d098bcdc75ea17448fd19d95477593/e2e/" + $newFileName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddrOverview, "", "", $newFileName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $cfgAddr, "", "", ".localization-config")

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Handoff never produced a usable target - clear the handoff file column
$wsZhCn.Range("C2").Clear()

$wsZhCn.Range("B2").Value = "Handoff transform failed"
$wsZhCn.Range("D2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Ignored"

$wsZhCn.Hyperlinks.Delete()
$mdAddrZhCn = "https://github.com/OpenLocalizationTest/oltest/blob/dbda4a81e1d098bcdc75ea17448fd19d95477593/e2e/" + $newFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddrZhCn, "", "", $newFileName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $cfgAddr, "", "", ".localization-config")

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Clear()

$wsDeDe.Range("B2").Value = "Handoff transform failed"
$wsDeDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Ignored"

$wsDeDe.Hyperlinks.Delete()
$mdAddrDeDe = "https://github.com/OpenLocalizationTest/oltest/blob/dbda4a81e1d098bcdc75ea17448fd19d95477593/e2e/" + $newFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddrDeDe, "", "", $newFileName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $cfgAddr, "", "", ".localization-config")
